# Release update: bump the embedded build timestamp from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet -------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Saranskaya Coal Mine, Kazakhstan, M1436, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ------------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 9; $row++) {
    $cell = $data.Range("S$row")
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
